$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$ws.Range("C2").Value = "MSG: None" + $nl + $nl + "MSG: The decision regarding Friday's movie has concluded without a selection." + $nl

$ws.Range("C3").Value = "MSG: None" + $nl + $nl + "MSG: The decision has been recorded successfully, and ""Barbie"" will be the movie shown on Friday." + $nl
$ws.Range("D3").Value = "Barbie_was_selected, "

$ws.Range("C4").Value = "MSG: None" + $nl + $nl + "MSG: The decision has been recorded as no decision about Friday's movie." + $nl

$ws.Range("C5").Value = "MSG: None" + $nl + $nl + "MSG: The rights to the movie ""Barbie"" have been successfully acquired for the upcoming Friday showing." + $nl

$ws.Range("C6").Value = "MSG: None" + $nl + $nl + "MSG: The decision has been recorded as no decision regarding the movie to be shown on Friday." + $nl

$ws.Range("C7").Value = "MSG: None" + $nl + $nl + "MSG: The movie rights for both ""Barbie"" and ""Oppenheimer"" have been successfully acquired for the assembly." + $nl
$ws.Range("D7").Value = "both_movies, "

$ws.Range("C8").Value = "MSG: None" + $nl + $nl + "MSG: The decision to acquire the rights for both movies has been recorded." + $nl
$ws.Range("D8").Value = "both_movies, "

$ws.Range("C9").Value = "MSG: None" + $nl + $nl + "MSG: I have successfully recorded the decision to acquire the rights for both movies." + $nl
